$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of data (row 7)
$ws.Range("A7").Value = "Vijay"
$ws.Range("B7").Value = "m"
$ws.Range("C7").Value = "n"
$ws.Range("D7").Value = 123
$ws.Range("E7").Value = "m@skrt.y"

# Add a mailto hyperlink on the email cell, matching the pattern of the
# existing rows (E2:E6), then reuse the same "Hyperlink" cell style that
# those rows already use.
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:m@skrt.y")
$ws.Range("E7").Style = $ws.Range("E6").Style

# Match the final selection left behind in the saved workbook.
$ws.Range("E12").Select()
